$wb = $excel.ActiveWorkbook

# ALC row 64 - Forged from the Void / Void Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5200
$ws.Range("I64").Value = 6000
$ws.Range("K64").Value = 6000
$ws.Range("M64").Value = -5752

# ALC row 67 - Dodging the Draft (L) / Void Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5200
$ws.Range("I67").Value = 6000
$ws.Range("K67").Value = 6000
$ws.Range("M67").Value = -5142

# ALC row 70 - Consecrating Congregation / Holy Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6617.8237
$ws.Range("I70").Value = 2225
$ws.Range("K70").Value = 6675
$ws.Range("M70").Value = -6405

# ALC row 73 - Curbing the Contagion (L) / Holy Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6617.8237
$ws.Range("I73").Value = 2225
$ws.Range("K73").Value = 6675
$ws.Range("M73").Value = -5739

# ALC row 96 - Scroll Down / Grade 1 Reisui of Intelligence
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2837.4
$ws.Range("I96").Value = 2325
$ws.Range("J96").Value = 4033
$ws.Range("K96").Value = 6975
$ws.Range("L96").Value = 12099
$ws.Range("M96").Value = -5602
$ws.Range("N96").Value = -14845

# ALC row 98 - The Dotted Line / Enchanted Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 557.2353000000001
$ws.Range("I98").Value = 557.2353000000001
$ws.Range("K98").Value = 557.2353000000001
$ws.Range("M98").Value = 940.7646999999999

# ALC row 112 - Making Ends Meet / Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2383.6924
$ws.Range("J112").Value = 2454.3635
$ws.Range("L112").Value = 7363.0905
$ws.Range("N112").Value = -9579.0905

# ALC row 122 - Wishful Inking / Enchanted High Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 557.2353000000001
$ws.Range("I122").Value = 557.2353000000001
$ws.Range("K122").Value = 1671.7059
$ws.Range("M122").Value = 778.2940999999998

# ALC row 137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2458.577
$ws.Range("I137").Value = 1041.9412
$ws.Range("K137").Value = 3125.8236
$ws.Range("M137").Value = -575.8235999999997

# ALC row 138 - All-night Crafting / Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2292.1365
$ws.Range("I138").Value = 1002.9091
$ws.Range("J138").Value = 2721.879
$ws.Range("K138").Value = 3008.7273
$ws.Range("L138").Value = 8165.637
$ws.Range("M138").Value = 2131.2727
$ws.Range("N138").Value = -18445.637

# ARM row 74 - As the Bolt Flies / Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2003.5312
$ws.Range("I74").Value = 1148.0385
$ws.Range("J74").Value = 5710.6665
$ws.Range("K74").Value = 1148.0385
$ws.Range("L74").Value = 5710.6665
$ws.Range("M74").Value = -274.0385000000001
$ws.Range("N74").Value = -7458.6665

# ARM row 77 - Heavy Metal Banned (L) / Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2003.5312
$ws.Range("I77").Value = 1148.0385
$ws.Range("J77").Value = 5710.6665
$ws.Range("K77").Value = 5740.192500000001
$ws.Range("L77").Value = 28553.3325
$ws.Range("M77").Value = -1372.192500000001
$ws.Range("N77").Value = -37289.3325

# ARM row 102 - Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1627.6875
$ws.Range("J102").Value = 1916.3334
$ws.Range("L102").Value = 1916.3334
$ws.Range("N102").Value = -5160.3334

# ARM row 122 - Haste for High Durium / High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 437546.7
$ws.Range("I122").Value = 1001107.5
$ws.Range("K122").Value = 3003322.5
$ws.Range("M122").Value = -3000872.5

# BSM row 20 - Smelt and Dealt / Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5531.636
$ws.Range("I20").Value = 6316.4443
$ws.Range("K20").Value = 6316.4443
$ws.Range("M20").Value = -6069.4443

# BSM row 94 - High Steal / High Steel Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1234.4375
$ws.Range("I94").Value = 1335.8572
$ws.Range("K94").Value = 1335.8572
$ws.Range("M94").Value = -884.8571999999999

# BSM row 105 - Ingot to Wing It / Molybdenum Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4912.1333
$ws.Range("I105").Value = 3390.25
$ws.Range("K105").Value = 3390.25
$ws.Range("M105").Value = -1643.25

# BSM row 134 - Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1904.3462
$ws.Range("I134").Value = 1574.5238
$ws.Range("J134").Value = 3289.6
$ws.Range("K134").Value = 4723.5714
$ws.Range("L134").Value = 9868.799999999999
$ws.Range("M134").Value = -2188.5714
$ws.Range("N134").Value = -14938.8

# CRP row 68 - Do You Even String Bow / Holy Cedar Composite Bow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 37999.75
$ws.Range("J68").Value = 37999.75
$ws.Range("L68").Value = 37999.75
$ws.Range("N68").Value = -39497.75

# CRP row 71 - Win One Bow, Get Three Free (L) / Holy Cedar Composite Bow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 37999.75
$ws.Range("J71").Value = 37999.75
$ws.Range("L71").Value = 113999.25
$ws.Range("N71").Value = -121487.25

# CRP row 99 - O Pine / Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13987.846
$ws.Range("J99").Value = 16585.715
$ws.Range("L99").Value = 16585.715
$ws.Range("N99").Value = -19581.715

# CRP row 126 - A Better Conductor / Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13987.846
$ws.Range("J126").Value = 16585.715
$ws.Range("L126").Value = 49757.145
$ws.Range("N126").Value = -54697.145

# CUL row 38 - Pretty as a Picture / Dark Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 111.07143
$ws.Range("I38").Value = 67.111115
$ws.Range("J38").Value = 190.2
$ws.Range("K38").Value = 201.333345
$ws.Range("L38").Value = 570.5999999999999
$ws.Range("M38").Value = 145.666655
$ws.Range("N38").Value = -1264.6

# CUL row 129 - Comfort Food / Yakow Moussaka
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1031.1666
$ws.Range("I129").Value = 630.8
$ws.Range("J129").Value = 3033
$ws.Range("K129").Value = 1892.4
$ws.Range("L129").Value = 9099
$ws.Range("M129").Value = 3107.6
$ws.Range("N129").Value = -19099

# GSM row 102 - Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3145.84
$ws.Range("I102").Value = 2163.8572
$ws.Range("K102").Value = 2163.8572
$ws.Range("M102").Value = -541.8571999999999

# GSM row 126 - Gold Rush Order / Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2952.4546
$ws.Range("I126").Value = 1265.75
$ws.Range("K126").Value = 3797.25
$ws.Range("M126").Value = -1327.25

# GSM row 132 - On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3636.08
$ws.Range("I132").Value = 2250
$ws.Range("J132").Value = 5137.6665
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 15412.9995
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -20472.9995

# LTW row 22 - Skin off Their Backs / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 18665.334
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 18665.334
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 18665.334
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -19255.334

# LTW row 27 - Fire and Hide / Aldgoat Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 18665.334
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 18665.334
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 18665.334
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -18879.334

# LTW row 61 - Spelling Me Softly / Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9109
$ws.Range("I61").Value = 11499
$ws.Range("J61").Value = 5922.3335
$ws.Range("K61").Value = 11499
$ws.Range("L61").Value = 5922.3335
$ws.Range("M61").Value = -11297
$ws.Range("N61").Value = -6326.3335

# LTW row 68 - You Could Say It's a Moving Target / Wyvern Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1999
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1250
$ws.Range("N68").ClearContents()

# LTW row 71 - They Call It Bloody Mary (L) / Wyvern Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9995
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6251
$ws.Range("N71").ClearContents()

# LTW row 100 - Tiger in the Sack / Tiger Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1796.4
$ws.Range("I100").Value = 1632.8889
$ws.Range("J100").Value = 2041.6666
$ws.Range("K100").Value = 1632.8889
$ws.Range("L100").Value = 2041.6666
$ws.Range("M100").Value = -1091.8889
$ws.Range("N100").Value = -3123.6666

# LTW row 113 - Peace in Rest / Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 9109
$ws.Range("I113").Value = 11499
$ws.Range("J113").Value = 5922.3335
$ws.Range("K113").Value = 11499
$ws.Range("L113").Value = 5922.3335
$ws.Range("M113").Value = -9329
$ws.Range("N113").Value = -10262.3335

# LTW row 122 - Hell on Leather / Gaja Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6763.174
$ws.Range("I122").Value = 7231.933
$ws.Range("K122").Value = 21695.799
$ws.Range("M122").Value = -19245.799

# LTW row 132 - Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5985.3076
$ws.Range("J132").Value = 6066.6665
$ws.Range("L132").Value = 18199.9995
$ws.Range("N132").Value = -23259.9995

# WVR row 62 - Pride Up in Smoke / Rainbow Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6370.5
$ws.Range("I62").Value = 4837.154
$ws.Range("J62").Value = 7699.4
$ws.Range("K62").Value = 4837.154
$ws.Range("L62").Value = 7699.4
$ws.Range("M62").Value = -4213.154
$ws.Range("N62").Value = -8947.4

# WVR row 65 - Desperate for Diversionaries (L) / Rainbow Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6370.5
$ws.Range("I65").Value = 4837.154
$ws.Range("J65").Value = 7699.4
$ws.Range("K65").Value = 24185.77
$ws.Range("L65").Value = 38497
$ws.Range("M65").Value = -21065.77
$ws.Range("N65").Value = -44737

# WVR row 81 - Where the Dragonflies, the Net Catches / Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3017.1333
$ws.Range("J81").Value = 2220.5
$ws.Range("L81").Value = 4441
$ws.Range("N81").Value = -6563

# WVR row 84 - To Kill a Dragon on Nameday (L) / Crawler Silk
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3017.1333
$ws.Range("J84").Value = 2220.5
$ws.Range("L84").Value = 22205
$ws.Range("N84").Value = -32813

# WVR row 122 - Heavy Armoire / Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2330.2
$ws.Range("I122").Value = 2255.7778
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6767.3334
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4317.3334
$ws.Range("N122").Value = -13900

# WVR row 136 - Weaving the Envelope / Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1535.625
$ws.Range("I136").Value = 1151.1945
$ws.Range("K136").Value = 3453.5835
$ws.Range("M136").Value = -903.5835000000002
